# Auto-generated script applying the Ixion_Profits market-data refresh
# (scheduled runner updates H/I/J/K/L/M/N market columns on affected rows)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4389750.5
$ws.Range("J17").Value = 4502304
$ws.Range("L17").Value = 13506912
$ws.Range("N17").Value = -13507248
$ws.Range("H62").Value = 1507.9286
$ws.Range("I62").Value = 1442.5834
$ws.Range("J62").Value = 1900
$ws.Range("K62").Value = 1442.5834
$ws.Range("L62").Value = 1900
$ws.Range("M62").Value = -818.5834
$ws.Range("N62").Value = -3148
$ws.Range("H65").Value = 1507.9286
$ws.Range("I65").Value = 1442.5834
$ws.Range("J65").Value = 1900
$ws.Range("K65").Value = 7212.916999999999
$ws.Range("L65").Value = 9500
$ws.Range("M65").Value = -4092.916999999999
$ws.Range("N65").Value = -15740
$ws.Range("H111").Value = 57267.5
$ws.Range("I111").Value = 1637.0769
$ws.Range("J111").Value = 201906.6
$ws.Range("K111").Value = 4911.2307
$ws.Range("L111").Value = 605719.8
$ws.Range("M111").Value = -1844.2307
$ws.Range("N111").Value = -611853.8
$ws.Range("H113").Value = 3324
$ws.Range("I113").Value = 3442
$ws.Range("J113").Value = 3176.5
$ws.Range("K113").Value = 3442
$ws.Range("L113").Value = 3176.5
$ws.Range("M113").Value = -188
$ws.Range("N113").Value = -9684.5
$ws.Range("H116").Value = 12960.5
$ws.Range("I116").Value = 18334.166
$ws.Range("J116").Value = 4900
$ws.Range("K116").Value = 18334.166
$ws.Range("L116").Value = 4900
$ws.Range("M116").Value = -14892.166
$ws.Range("N116").Value = -11784
$ws.Range("H137").Value = 1496.6129
$ws.Range("I137").Value = 1343.8
$ws.Range("J137").Value = 2133.3333
$ws.Range("K137").Value = 4031.4
$ws.Range("L137").Value = 6399.999899999999
$ws.Range("M137").Value = -1481.4
$ws.Range("N137").Value = -11499.9999
$ws.Range("H138").Value = 2248.0833
$ws.Range("I138").Value = 700.61365
$ws.Range("J138").Value = 3950.3
$ws.Range("K138").Value = 2101.84095
$ws.Range("L138").Value = 11850.9
$ws.Range("M138").Value = 3038.15905
$ws.Range("N138").Value = -22130.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8167.357
$ws.Range("I45").Value = 10104.728
$ws.Range("J45").Value = 1063.6666
$ws.Range("K45").Value = 10104.728
$ws.Range("L45").Value = 1063.6666
$ws.Range("M45").Value = -9727.727999999999
$ws.Range("N45").Value = -1817.6666
$ws.Range("H74").Value = 1245.7354
$ws.Range("I74").Value = 1090.9259
$ws.Range("K74").Value = 1090.9259
$ws.Range("M74").Value = -216.9259
$ws.Range("H77").Value = 1245.7354
$ws.Range("I77").Value = 1090.9259
$ws.Range("K77").Value = 5454.6295
$ws.Range("M77").Value = -1086.6295

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2482511.8
$ws.Range("I16").Value = 4274322.5
$ws.Range("J16").Value = 1543.3077
$ws.Range("K16").Value = 4274322.5
$ws.Range("L16").Value = 1543.3077
$ws.Range("M16").Value = -4274035.5
$ws.Range("N16").Value = -2117.3077
$ws.Range("H31").Value = 4444.1875
$ws.Range("I31").Value = 2224.889
$ws.Range("J31").Value = 7297.5713
$ws.Range("K31").Value = 2224.889
$ws.Range("L31").Value = 7297.5713
$ws.Range("M31").Value = -1929.889
$ws.Range("N31").Value = -7887.5713
$ws.Range("H34").Value = 4444.1875
$ws.Range("I34").Value = 2224.889
$ws.Range("J34").Value = 7297.5713
$ws.Range("K34").Value = 2224.889
$ws.Range("L34").Value = 7297.5713
$ws.Range("M34").Value = -2022.889
$ws.Range("N34").Value = -7701.5713
$ws.Range("H94").Value = 3232.1072
$ws.Range("I94").Value = 3939.6
$ws.Range("K94").Value = 3939.6
$ws.Range("M94").Value = -3488.6
$ws.Range("H113").Value = 2482511.8
$ws.Range("I113").Value = 4274322.5
$ws.Range("J113").Value = 1543.3077
$ws.Range("K113").Value = 4274322.5
$ws.Range("L113").Value = 1543.3077
$ws.Range("M113").Value = -4272152.5
$ws.Range("N113").Value = -5883.3077

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 400
$ws.Range("I31").Value = 400
$ws.Range("K31").Value = 1200
$ws.Range("M31").Value = -912
$ws.Range("H132").Value = 1575
$ws.Range("J132").Value = 2533.3333
$ws.Range("L132").Value = 22799.9997
$ws.Range("N132").Value = -27859.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 207121010
$ws.Range("J11").Value = 9201669
$ws.Range("L11").Value = 9201669
$ws.Range("N11").Value = -9201947
$ws.Range("H122").Value = 1411056.8
$ws.Range("I122").Value = 2236693.8
$ws.Range("J122").Value = 2617.4119
$ws.Range("K122").Value = 6710081.399999999
$ws.Range("L122").Value = 7852.2357
$ws.Range("M122").Value = -6707631.399999999
$ws.Range("N122").Value = -12752.2357
$ws.Range("H126").Value = 6961.75
$ws.Range("I126").Value = 8502.267
$ws.Range("J126").Value = 2340.2
$ws.Range("K126").Value = 25506.801
$ws.Range("L126").Value = 7020.599999999999
$ws.Range("M126").Value = -23036.801
$ws.Range("N126").Value = -11960.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 8207.143
$ws.Range("I29").Value = 7900
$ws.Range("J29").Value = 8258.333000000001
$ws.Range("K29").Value = 7900
$ws.Range("L29").Value = 8258.333000000001
$ws.Range("M29").Value = -7605
$ws.Range("N29").Value = -8848.333000000001
$ws.Range("H43").Value = 8510
$ws.Range("J43").Value = 8510
$ws.Range("L43").Value = 8510
$ws.Range("N43").Value = -8896
$ws.Range("H68").Value = 71431680
$ws.Range("I68").Value = 3209
$ws.Range("J68").Value = 500002500
$ws.Range("K68").Value = 3209
$ws.Range("L68").Value = 500002500
$ws.Range("M68").Value = -2460
$ws.Range("N68").Value = -500003998
$ws.Range("H71").Value = 71431680
$ws.Range("I71").Value = 3209
$ws.Range("J71").Value = 500002500
$ws.Range("K71").Value = 16045
$ws.Range("L71").Value = 2500012500
$ws.Range("M71").Value = -12301
$ws.Range("N71").Value = -2500019988
$ws.Range("H93").Value = 62500868
$ws.Range("I93").Value = 992.8570999999999
$ws.Range("J93").Value = 500000000
$ws.Range("K93").Value = 992.8570999999999
$ws.Range("L93").Value = 500000000
$ws.Range("M93").Value = 255.1429000000001
$ws.Range("N93").Value = -500002496
$ws.Range("H100").Value = 1358.25
$ws.Range("I100").Value = 1377.2222
$ws.Range("K100").Value = 1377.2222
$ws.Range("M100").Value = -836.2221999999999
$ws.Range("H132").Value = 9168803
$ws.Range("I132").Value = 12389327
$ws.Range("J132").Value = 2694.6924
$ws.Range("K132").Value = 37167981
$ws.Range("L132").Value = 8084.0772
$ws.Range("M132").Value = -37165451
$ws.Range("N132").Value = -13144.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 6500
$ws.Range("J32").Value = 6500
$ws.Range("L32").Value = 6500
$ws.Range("N32").Value = -7134
$ws.Range("H34").Value = 6500
$ws.Range("J34").Value = 6500
$ws.Range("L34").Value = 6500
$ws.Range("N34").Value = -6906
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184
$ws.Range("H136").Value = 966
$ws.Range("I136").Value = 496.73334
$ws.Range("J136").Value = 1971.5714
$ws.Range("K136").Value = 1490.20002
$ws.Range("L136").Value = 5914.7142
$ws.Range("M136").Value = 1059.79998
$ws.Range("N136").Value = -11014.7142
